# Scheduled refresh: update market-price derived columns (H-N) on several
# Leve-profit sheets. Values below come straight from the upstream data pull;
# no formulas are involved (columns are plain numeric snapshots).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 4
$ws.Range("H4").Value = 500
$ws.Range("I4").Value = 132
$ws.Range("J4").Value = 1481.3334
$ws.Range("K4").Value = 132
$ws.Range("L4").Value = 1481.3334
$ws.Range("M4").Value = -18
$ws.Range("N4").Value = -1709.3334

# Row 7
$ws.Range("H7").Value = 999.25
$ws.Range("J7").Value = 999.25
$ws.Range("L7").Value = 999.25
$ws.Range("N7").Value = -1223.25

# Row 14
$ws.Range("H14").Value = 999.25
$ws.Range("J14").Value = 999.25
$ws.Range("L14").Value = 999.25
$ws.Range("N14").Value = -1381.25

# Row 64
$ws.Range("H64").Value = 22225736
$ws.Range("I64").Value = 33335548
$ws.Range("J64").Value = 6109.8
$ws.Range("K64").Value = 33335548
$ws.Range("L64").Value = 6109.8
$ws.Range("M64").Value = -33335300
$ws.Range("N64").Value = -6605.8

# Row 67
$ws.Range("H67").Value = 22225736
$ws.Range("I67").Value = 33335548
$ws.Range("J67").Value = 6109.8
$ws.Range("K67").Value = 33335548
$ws.Range("L67").Value = 6109.8
$ws.Range("M67").Value = -33334690
$ws.Range("N67").Value = -7825.8

# Row 70
$ws.Range("H70").Value = 8772.556
$ws.Range("J70").Value = 8772.556
$ws.Range("L70").Value = 26317.668
$ws.Range("N70").Value = -26857.668

# Row 73
$ws.Range("H73").Value = 8772.556
$ws.Range("J73").Value = 8772.556
$ws.Range("L73").Value = 26317.668
$ws.Range("N73").Value = -28189.668

# Row 86
$ws.Range("H86").Value = 1490.3846
$ws.Range("I86").Value = 1649.7778
$ws.Range("K86").Value = 1649.7778
$ws.Range("M86").Value = -526.7778000000001

# Row 89
$ws.Range("H89").Value = 1490.3846
$ws.Range("I89").Value = 1649.7778
$ws.Range("K89").Value = 8248.889000000001
$ws.Range("M89").Value = -2632.889000000001

# Row 138
$ws.Range("H138").Value = 3144.5095
$ws.Range("I138").Value = 1450.2258
$ws.Range("J138").Value = 5531.909
$ws.Range("K138").Value = 4350.6774
$ws.Range("L138").Value = 16595.727
$ws.Range("M138").Value = 789.3226000000004
$ws.Range("N138").Value = -26875.727

$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Range("H61").Value = 4469.271
$ws.Range("I61").Value = 4272.5713
$ws.Range("J61").Value = 4577.9736
$ws.Range("K61").Value = 4272.5713
$ws.Range("L61").Value = 4577.9736
$ws.Range("M61").Value = -4060.5713
$ws.Range("N61").Value = -5001.9736

# Row 132
$ws.Range("H132").Value = 533151.1
$ws.Range("I132").Value = 535293.7
$ws.Range("K132").Value = 1605881.1
$ws.Range("M132").Value = -1603351.1

# Row 136
$ws.Range("H136").Value = 4469.271
$ws.Range("I136").Value = 4272.5713
$ws.Range("J136").Value = 4577.9736
$ws.Range("K136").Value = 12817.7139
$ws.Range("L136").Value = 13733.9208
$ws.Range("M136").Value = -10267.7139
$ws.Range("N136").Value = -18833.9208

$ws = $wb.Worksheets.Item("BSM")
# Row 22
$ws.Range("H22").Value = 1949
$ws.Range("I22").Value = 337
$ws.Range("J22").Value = 4205.8
$ws.Range("K22").Value = 337
$ws.Range("L22").Value = 4205.8
$ws.Range("M22").Value = -164
$ws.Range("N22").Value = -4551.8

# Row 134
$ws.Range("H134").Value = 985284.25
$ws.Range("I134").Value = 1140575.4
$ws.Range("J134").Value = 9168.143
$ws.Range("K134").Value = 3421726.2
$ws.Range("L134").Value = 27504.429
$ws.Range("M134").Value = -3419191.2
$ws.Range("N134").Value = -32574.429

$ws = $wb.Worksheets.Item("CRP")
# Row 58
$ws.Range("H58").Value = 33345220
$ws.Range("I58").Value = 40008664
$ws.Range("K58").Value = 40008664
$ws.Range("M58").Value = -40008461

# Row 132
$ws.Range("H132").Value = 14632.3
$ws.Range("J132").Value = 31923
$ws.Range("L132").Value = 95769
$ws.Range("N132").Value = -100829

# Row 134
$ws.Range("H134").Value = 62510372
$ws.Range("I134").Value = 83339270
$ws.Range("J134").Value = 23662.5
$ws.Range("K134").Value = 250017810
$ws.Range("L134").Value = 70987.5
$ws.Range("M134").Value = -250015275
$ws.Range("N134").Value = -76057.5

# Row 136
$ws.Range("H136").Value = 33345220
$ws.Range("I136").Value = 40008664
$ws.Range("K136").Value = 120025992
$ws.Range("M136").Value = -120023442

$ws = $wb.Worksheets.Item("CUL")
# Row 3
$ws.Range("H3").Value = 13997.5
$ws.Range("I3").Value = 2601.4285
$ws.Range("K3").Value = 7804.2855
$ws.Range("M3").Value = -7692.2855

# Row 16
$ws.Range("H16").Value = 0
$ws.Range("J16").Value = 0
$ws.Range("L16").Value = 0
$ws.Range("N16").ClearContents()

# Row 37
$ws.Range("H37").Value = 103727.836
$ws.Range("J37").Value = 103727.836
$ws.Range("L37").Value = 311183.508
$ws.Range("N37").Value = -311407.508

# Row 40
$ws.Range("H40").Value = 77.94444
$ws.Range("I40").Value = 83.916664
$ws.Range("J40").Value = 66
$ws.Range("K40").Value = 335.666656
$ws.Range("L40").Value = 264
$ws.Range("M40").Value = -266.666656
$ws.Range("N40").Value = -402

# Row 113
$ws.Range("H113").Value = 3473013.8
$ws.Range("I113").Value = 10417020
$ws.Range("J113").Value = 1010.5833
$ws.Range("K113").Value = 31251060
$ws.Range("L113").Value = 3031.7499
$ws.Range("M113").Value = -31248890
$ws.Range("N113").Value = -7371.7499

# Row 131
$ws.Range("H131").Value = 37683900
$ws.Range("I131").Value = 28072542
$ws.Range("J131").Value = 83337864
$ws.Range("K131").Value = 84217626
$ws.Range("L131").Value = 250013592
$ws.Range("M131").Value = -84212586
$ws.Range("N131").Value = -250023672

$ws = $wb.Worksheets.Item("GSM")
# Row 2
$ws.Range("H2").Value = 61.666668
$ws.Range("J2").Value = 78.28570999999999
$ws.Range("L2").Value = 78.28570999999999
$ws.Range("N2").Value = -304.28571

# Row 107
$ws.Range("H107").Value = 538.55554
$ws.Range("I107").Value = 160.42857
$ws.Range("J107").Value = 779.1818
$ws.Range("K107").Value = 160.42857
$ws.Range("L107").Value = 779.1818
$ws.Range("M107").Value = 1759.57143
$ws.Range("N107").Value = -4619.1818

# Row 132
$ws.Range("H132").Value = 5482.486
$ws.Range("I132").Value = 5555.5586
$ws.Range("J132").Value = 2998
$ws.Range("K132").Value = 16666.6758
$ws.Range("L132").Value = 8994
$ws.Range("M132").Value = -14136.6758
$ws.Range("N132").Value = -14054

# Row 138
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 10689.4
$ws.Range("J7").Value = 16032.333
$ws.Range("L7").Value = 16032.333
$ws.Range("N7").Value = -16256.333

# Row 22
$ws.Range("H22").Value = 976.4286
$ws.Range("I22").Value = 968.75
$ws.Range("K22").Value = 968.75
$ws.Range("M22").Value = -673.75

# Row 27
$ws.Range("H27").Value = 976.4286
$ws.Range("I27").Value = 968.75
$ws.Range("K27").Value = 968.75
$ws.Range("M27").Value = -861.75

# Row 40
$ws.Range("H40").Value = 20576.842
$ws.Range("I40").Value = 27767.46
$ws.Range("K40").Value = 27767.46
$ws.Range("M40").Value = -27631.46

# Row 46
$ws.Range("H46").Value = 1910.579
$ws.Range("I46").Value = 1086.4286
$ws.Range("J46").Value = 2391.3333
$ws.Range("K46").Value = 1086.4286
$ws.Range("L46").Value = 2391.3333
$ws.Range("M46").Value = -898.4286
$ws.Range("N46").Value = -2767.3333

# Row 126
$ws.Range("H126").Value = 10689.4
$ws.Range("J126").Value = 16032.333
$ws.Range("L126").Value = 48096.999
$ws.Range("N126").Value = -53036.999

# Row 132
$ws.Range("H132").Value = 2505.4285
$ws.Range("I132").Value = 2590
$ws.Range("J132").Value = 1998
$ws.Range("K132").Value = 7770
$ws.Range("L132").Value = 5994
$ws.Range("M132").Value = -5240
$ws.Range("N132").Value = -11054

$ws = $wb.Worksheets.Item("WVR")
# Row 122
$ws.Range("H122").Value = 29220.666
$ws.Range("I122").Value = 25997.666
$ws.Range("J122").Value = 30832.166
$ws.Range("K122").Value = 77992.99800000001
$ws.Range("L122").Value = 92496.49800000001
$ws.Range("M122").Value = -75542.99800000001
$ws.Range("N122").Value = -97396.49800000001

# Row 132
$ws.Range("H132").Value = 12362
$ws.Range("I132").Value = 7747.95
$ws.Range("K132").Value = 23243.85
$ws.Range("M132").Value = -20713.85

# Row 136
$ws.Range("H136").Value = 16679545
$ws.Range("I136").Value = 19240110
$ws.Range("J136").Value = 35874.25
$ws.Range("K136").Value = 57720330
$ws.Range("L136").Value = 107622.75
$ws.Range("M136").Value = -57717780
$ws.Range("N136").Value = -112722.75

